$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the new values for row 14 (record #13)
$ws.Range("B14").Value = "Al-azazi & Ghurab 2022"
$ws.Range("E14").Value = "ANN-LSTM"

# Match formatting of neighbouring data cells in the same row (e.g. A14): wrap text, no special fill/font
$ws.Range("B14").WrapText = $true
$ws.Range("E14").WrapText = $true

# Row grows to 2 lines of wrapped text -> 30pt tall, matching the other wrapped rows
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(14).RowHeight = 30

# Update the frozen pane top-left cell and the active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("F14").Select()
